# Scheduled-runner update: refresh market-board derived profit figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H-N)
# across several leve rows on the ALC, ARM, BSM, CRP, GSM, LTW and WVR
# sheets of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1843.4
$ws.Range("J17").Value = 1843.4
$ws.Range("L17").Value = 5530.200000000001
$ws.Range("N17").Value = -5866.200000000001

$ws.Range("H41").Value = 1841.9286
$ws.Range("I41").Value = 3019.625
$ws.Range("J41").Value = 271.66666
$ws.Range("K41").Value = 3019.625
$ws.Range("L41").Value = 271.66666
$ws.Range("M41").Value = -2579.625
$ws.Range("N41").Value = -1151.66666

$ws.Range("H64").Value = 6250.5
$ws.Range("I64").Value = 5002
$ws.Range("K64").Value = 5002
$ws.Range("M64").Value = -4754

$ws.Range("H67").Value = 6250.5
$ws.Range("I67").Value = 5002
$ws.Range("K67").Value = 5002
$ws.Range("M67").Value = -4144

$ws.Range("H70").Value = 62502210
$ws.Range("J70").Value = 62502210
$ws.Range("L70").Value = 187506630
$ws.Range("N70").Value = -187507170

$ws.Range("H73").Value = 62502210
$ws.Range("J73").Value = 62502210
$ws.Range("L73").Value = 187506630
$ws.Range("N73").Value = -187508502

$ws.Range("H74").Value = 6249.25
$ws.Range("I74").Value = 6249.25
$ws.Range("K74").Value = 6249.25
$ws.Range("M74").Value = -5313.25

$ws.Range("H77").Value = 6249.25
$ws.Range("I77").Value = 6249.25
$ws.Range("K77").Value = 31246.25
$ws.Range("M77").Value = -26566.25

$ws.Range("H98").Value = 1418.091
$ws.Range("I98").Value = 1491.5
$ws.Range("J98").Value = 953.1667
$ws.Range("K98").Value = 1491.5
$ws.Range("L98").Value = 953.1667
$ws.Range("M98").Value = 6.5
$ws.Range("N98").Value = -3949.1667

$ws.Range("H122").Value = 1418.091
$ws.Range("I122").Value = 1491.5
$ws.Range("J122").Value = 953.1667
$ws.Range("K122").Value = 4474.5
$ws.Range("L122").Value = 2859.5001
$ws.Range("M122").Value = -2024.5
$ws.Range("N122").Value = -7759.5001

$ws.Range("H125").Value = 1849.5
$ws.Range("I125").Value = 1849
$ws.Range("K125").Value = 16641
$ws.Range("M125").Value = -14181

$ws.Range("H132").Value = 3638.3062
$ws.Range("I132").Value = 4201.343
$ws.Range("K132").Value = 12604.029
$ws.Range("M132").Value = -10074.029

$ws.Range("H137").Value = 32490.408
$ws.Range("I137").Value = 40592.715
$ws.Range("J137").Value = 4132.3335
$ws.Range("K137").Value = 121778.145
$ws.Range("L137").Value = 12397.0005
$ws.Range("M137").Value = -119228.145
$ws.Range("N137").Value = -17497.0005

$ws.Range("H138").Value = 2811.8096
$ws.Range("J138").Value = 3649.8333
$ws.Range("L138").Value = 10949.4999
$ws.Range("N138").Value = -21229.4999

$ws.Range("H141").Value = 1889.3636
$ws.Range("I141").Value = 1288.45
$ws.Range("J141").Value = 7898.5
$ws.Range("K141").Value = 3865.35
$ws.Range("L141").Value = 23695.5
$ws.Range("M141").Value = 1314.65
$ws.Range("N141").Value = -34055.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3055.2273
$ws.Range("I61").Value = 2721.842
$ws.Range("K61").Value = 2721.842
$ws.Range("M61").Value = -2509.842

$ws.Range("H63").Value = 3779.5
$ws.Range("I63").Value = 2333.3333
$ws.Range("J63").Value = 4399.2856
$ws.Range("K63").Value = 2333.3333
$ws.Range("L63").Value = 4399.2856
$ws.Range("M63").Value = -1647.3333
$ws.Range("N63").Value = -5771.2856

$ws.Range("H66").Value = 3779.5
$ws.Range("I66").Value = 2333.3333
$ws.Range("J66").Value = 4399.2856
$ws.Range("K66").Value = 11666.6665
$ws.Range("L66").Value = 21996.428
$ws.Range("M66").Value = -8234.6665
$ws.Range("N66").Value = -28860.428

$ws.Range("H132").Value = 2912.739
$ws.Range("I132").Value = 2420.2354
$ws.Range("K132").Value = 7260.706200000001
$ws.Range("M132").Value = -4730.706200000001

$ws.Range("H136").Value = 3055.2273
$ws.Range("I136").Value = 2721.842
$ws.Range("K136").Value = 8165.526
$ws.Range("M136").Value = -5615.526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15877161
$ws.Range("I134").Value = 1498.9166
$ws.Range("J134").Value = 37044708
$ws.Range("K134").Value = 4496.7498
$ws.Range("L134").Value = 111134124
$ws.Range("M134").Value = -1961.7498
$ws.Range("N134").Value = -111139194

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1488.4706
$ws.Range("I132").Value = 1193.6538
$ws.Range("K132").Value = 3580.9614
$ws.Range("M132").Value = -1050.9614

$ws.Range("H134").Value = 1744.2106
$ws.Range("I134").Value = 1281.2727
$ws.Range("K134").Value = 3843.8181
$ws.Range("M134").Value = -1308.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19770.053
$ws.Range("I70").Value = 102525
$ws.Range("J70").Value = 4253.5
$ws.Range("K70").Value = 102525
$ws.Range("L70").Value = 4253.5
$ws.Range("M70").Value = -102255
$ws.Range("N70").Value = -4793.5

$ws.Range("H73").Value = 19770.053
$ws.Range("I73").Value = 102525
$ws.Range("J73").Value = 4253.5
$ws.Range("K73").Value = 102525
$ws.Range("L73").Value = 4253.5
$ws.Range("M73").Value = -101589
$ws.Range("N73").Value = -6125.5

$ws.Range("H102").Value = 2958.8235
$ws.Range("I102").Value = 2066.6667
$ws.Range("K102").Value = 2066.6667
$ws.Range("M102").Value = -444.6667000000002

$ws.Range("H122").Value = 2716.3845
$ws.Range("I122").Value = 2164.3
$ws.Range("K122").Value = 6492.900000000001
$ws.Range("M122").Value = -4042.900000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3566.1428
$ws.Range("I46").Value = 698.75
$ws.Range("K46").Value = 698.75
$ws.Range("M46").Value = -510.75

$ws.Range("H82").Value = 3755.5
$ws.Range("I82").Value = 3737.25
$ws.Range("J82").Value = 3762.8
$ws.Range("K82").Value = 3737.25
$ws.Range("L82").Value = 3762.8
$ws.Range("M82").Value = -3376.25
$ws.Range("N82").Value = -4484.8

$ws.Range("H85").Value = 3755.5
$ws.Range("I85").Value = 3737.25
$ws.Range("J85").Value = 3762.8
$ws.Range("K85").Value = 3737.25
$ws.Range("L85").Value = 3762.8
$ws.Range("M85").Value = -2489.25
$ws.Range("N85").Value = -6258.8

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H132").Value = 3560.818
$ws.Range("I132").Value = 3018.7778
$ws.Range("K132").Value = 9056.3334
$ws.Range("M132").Value = -6526.3334

$ws.Range("H136").Value = 3893.8948
$ws.Range("I136").Value = 3835.6667
$ws.Range("K136").Value = 11507.0001
$ws.Range("M136").Value = -8957.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 55187.5
$ws.Range("I51").Value = 50000
$ws.Range("J51").Value = 56916.668
$ws.Range("K51").Value = 50000
$ws.Range("L51").Value = 56916.668
$ws.Range("M51").Value = -49490
$ws.Range("N51").Value = -57936.668

$ws.Range("H113").Value = 300
$ws.Range("I113").Value = 242.73334
$ws.Range("K113").Value = 728.20002
$ws.Range("M113").Value = 1441.79998

$ws.Range("H125").Value = 99816.336
$ws.Range("J125").Value = 99816.336
$ws.Range("L125").Value = 99816.336
$ws.Range("N125").Value = -109656.336

$ws.Range("H132").Value = 2370.2712
$ws.Range("I132").Value = 2001.7925
$ws.Range("K132").Value = 6005.377500000001
$ws.Range("M132").Value = -3475.377500000001

$ws.Range("H136").Value = 32189.484
$ws.Range("I136").Value = 889.73914
$ws.Range("K136").Value = 2669.21742
$ws.Range("M136").Value = -119.2174199999999
